# CIERRE 13 Ago 22
# Move the payroll sheet from "SEMANA 31 (Aug 1-7, 2022)" to
# "SEMANA 32 (Aug 8-14, 2022)": update the week-label cell, zero out the
# one-off "EXTRAS" amount for the new week, and move the selection.
# All other formulas (TODAY()-based dates, the SEMANA references on rows
# 27/43, and the E41 total) recompute automatically on recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("recibos")

# Week label - only B9 holds the literal text; H9/B27/H27/B43 all derive
# from it via formulas (=B9, =B9, =B27, =H27) and recalc on their own.
$ws.Range("B9").Value = "SEMANA  32  DEL    08      Al   14   DE   AGOSTO          2022"

# "EXTRAS" amount for this period goes to 0 (was 1250); E41 = SUM(E38:E40)
# recalculates to 2500 automatically.
$ws.Range("E40").Value = 0

# Move the on-screen selection to I37:I38 (active cell I38).
$ws.Activate()
$ws.Range("I37:I38").Select()
